# Needle calibration data: sort the data rows (A2:D12) ascending by
# column A (time) now that the needle calibration has been performed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D12")
$values = $dataRange.Value2

$nRows = 11
$nCols = 4

# Row order (0-based, relative to the range) to be sorted ascending by
# column A (time).
$order = @(0,1,2,3,4,5,6,7,8,9,10)

# Simple bubble sort of $order using column A as the key.
for ($i = 0; $i -lt $nRows - 1; $i++) {
    for ($j = 0; $j -lt $nRows - 1 - $i; $j++) {
        $a = $values[$order[$j] + 1, 1]
        $b = $values[$order[$j + 1] + 1, 1]
        if ($a -gt $b) {
            $tmp = $order[$j]
            $order[$j] = $order[$j + 1]
            $order[$j + 1] = $tmp
        }
    }
}

$rowOut = 2
foreach ($srcRow in $order) {
    for ($c = 1; $c -le $nCols; $c++) {
        $ws.Cells.Item($rowOut, $c).Value = $values[$srcRow + 1, $c]
    }
    $rowOut = $rowOut + 1
}
